# Ajout au DF de donnees geographiques et politiques
# (ancien ministre/membre du gouvernement)
#
# - la colonne "depart_CORSE" est supprimee (les infos utiles qu'elle
#   portait n'existent plus : on se base desormais sur depart_frontalier)
# - la colonne "depart_OM" est renommee "depart_DOM" et etendue aux
#   collectivites d'outre-mer qui n'y figuraient pas encore
# - la Corse (Corse-du-Sud / Haute-Corse) est desormais marquee comme
#   "depart_frontalier"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Supprime la colonne "depart_CORSE" (colonne E) : "region_parisienne"
#    (colonne F) se decale donc en colonne E.
$ws.Range("E:E").Delete()

# 2) Renomme l'ancienne colonne "depart_OM" (colonne D) en "depart_DOM".
$ws.Cells.Item(1, 4).Value2 = "depart_DOM"

# 3) Marque les collectivites d'outre-mer restantes comme "depart_DOM".
$domAjouts = @(
    "LA-REUNION",
    "SAINT-PIERRE-ET-MIQUELON",
    "SAINT-MARTIN/SAINT-BARTHELEMY",
    "WALLIS-ET-FUTUNA",
    "POLYNESIE-FRANCAISE",
    "NOUVELLE-CALEDONIE"
)
foreach ($dep in $domAjouts) {
    $cell = $ws.Columns.Item(2).Find($dep)
    $ws.Cells.Item($cell.Row, 4).Value2 = 1
}

# 4) La Corse (Corse-du-Sud, Haute-Corse) passe a "depart_frontalier" = 1.
$corseDeps = @("CORSE-DU-SUD", "HAUTE-CORSE")
foreach ($dep in $corseDeps) {
    $cell = $ws.Columns.Item(2).Find($dep)
    $ws.Cells.Item($cell.Row, 3).Value2 = 1
}

# 5) Replique la selection active enregistree dans le classeur d'origine.
$ws.Range("C2").Select() | Out-Null
